$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1 and Q1, matching style of O1 (bold/border/centered)
foreach ($addr in @("P1", "Q1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
    $c.Borders.Weight = 2
}
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing columns I, K, M, O for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 15).Value = 1   # O
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
